$wb = $excel.ActiveWorkbook

# Rename sheets: "add_spe_" -> "add_pe_", "add_s_" -> "" (dropped prefix),
# and "add_s_SAP_Sinis_Ced" -> "SAP_Sinis_Ced"
$wb.Worksheets.Item("add_spe_Canal-Poliza").Name = "add_pe_Canal-Poliza"
$wb.Worksheets.Item("add_spe_Canal-Canal").Name = "add_pe_Canal-Canal"
$wb.Worksheets.Item("add_spe_Canal-Sucursal").Name = "add_pe_Canal-Sucursal"
$wb.Worksheets.Item("add_spe_Amparos").Name = "add_pe_Amparos"
$wb.Worksheets.Item("add_s_Atipicos").Name = "Atipicos"
$wb.Worksheets.Item("add_s_Inc_Ced_Atipicos").Name = "Inc_Ced_Atipicos"
$wb.Worksheets.Item("add_s_SAP_Sinis_Ced").Name = "SAP_Sinis_Ced"

# Update the selection on the active sheet ("SAP_Sinis_Ced") from I13 to F15
$ws = $wb.Worksheets.Item("SAP_Sinis_Ced")
$ws.Activate()
$ws.Range("F15").Select()
